$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing
# "fbref_team_names" (A) -> B and "vaibhav_team_names" (B) -> C.
$ws.Range("A1").EntireColumn.Insert()

# New column D: vaibhav_master_team_id (numeric row index 0..26)
$ws.Range("D1").Value = "vaibhav_master_team_id"
for ($i = 0; $i -le 26; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $i
}

# New column A: fbref_master_team_id (numeric row index 0..26)
$ws.Range("A1").Value = "fbref_master_team_id"
for ($i = 0; $i -le 26; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}

# Update the visible window/selection to match the saved view state.
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("F8").Select()
